# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") computed from the regenerated save data.
$kValues = @{
    2  = 1
    3  = 0
    4  = 4
    5  = 0
    6  = 1
    7  = 3
    8  = 2
    9  = 1
    10 = 3
    11 = 1
    12 = 2
    13 = 3
    14 = 1
    15 = 0
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 2
    22 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
